$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.372.72"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "'1.873.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("D4").Value = "'1.0000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'238.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.09%  "
$ws.Range("D6").Value = "'0.9999"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("E7").Value = "  -0.55%  "
$ws.Range("D8").Value = "'0.2825"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.69%  "
$ws.Range("D9").Value = "'0.06530"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.83%  "
$ws.Range("D10").Value = "'1.872.76"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.39%  "
$ws.Range("D11").Value = "'0.07481"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.61%  "
$ws.Range("D12").Value = "'16.55"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.92%  "
$ws.Range("D13").Value = "'5.096"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.52%  "
$ws.Range("D14").Value = "'88.20"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.26%  "
$ws.Range("D15").Value = "'0.6585"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.52%  "
$ws.Range("D16").Value = "'30.358.20"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("D17").Value = "'13.30"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("D18").Value = "'0.9999"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'0.000007629"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.29%  "
$ws.Range("D20").Value = "'2.118.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "'224.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +15.25%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'5.307"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.34%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").Value = "'6.191"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.18%  "
$ws.Range("D25").Value = "'9.284"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.12%  "
$ws.Range("D26").Value = "'165.92"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.73%  "
$ws.Range("D27").Value = "'18.81"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.09%  "
$ws.Range("D28").Value = "'1.985"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.67%  "
$ws.Range("E29").Value = "  +1.40%  "
$ws.Range("D30").Value = "'0.09399"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.12%  "
$ws.Range("D32").Value = "'4.019"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.71%  "
$ws.Range("D33").Value = "'0.05055"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.70%  "
$ws.Range("D34").Value = "'1.227"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +12.13%  "
$ws.Range("D35").Value = "'0.7524"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.88%  "
$ws.Range("D36").Value = "'2.699"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.25%  "
$ws.Range("D37").Value = "'0.01840"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.60%  "
$ws.Range("D38").Value = "'2.618"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.79%  "
$ws.Range("D39").Value = "'2.090"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.51%  "
$ws.Range("D40").Value = "'0.9073"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.35%  "
$ws.Range("D41").Value = "'5.934"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.54%  "
$ws.Range("D42").Value = "'106.99"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.71%  "
$ws.Range("D43").Value = "'0.4304"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.79%  "
$ws.Range("E44").Value = "  +0.33%  "
$ws.Range("D45").Value = "'7.462"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.09%  "
$ws.Range("E46").Value = "  -0.99%  "
$ws.Range("D47").Value = "'64.88"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.29%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'9.105"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.21%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'1.495"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +9.86%  "
$ws.Range("D50").Value = "'34.20"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.21%  "
$ws.Range("D51").Value = "'0.3926"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.94%  "
